$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 158 ("弟が歩くぞ" post), shifting subsequent rows up.
$ws.Rows.Item(158).Delete()
